$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values change
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 data changes - B2 and D2 cleared, C2 and E2 updated
$ws.Range("B2").ClearContents()
$ws.Range("C2").Value = -1.5036895116748712
$ws.Range("D2").ClearContents()
$ws.Range("E2").Value = -1.4615611087034994

# Row 3 data changes
$ws.Range("B3").Value = -2.087878713081964
$ws.Range("C3").Value = -0.71566236437726283
$ws.Range("D3").Value = -2.561795311586474
$ws.Range("E3").Value = 1.6914956341777456

# Update selection to match new sqref B1:E3
$ws.Range("B1:E3").Select()
